$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.004", "1.000") are preserved verbatim as text, matching
# the inlineStr/text semantics of the original cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.235.08"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "1.792.32"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "338.01"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").Value = "0.4525"
$ws.Range("E7").Value = "  +20.12%  "

$ws.Range("D8").Value = "0.3572"
$ws.Range("E8").Value = "  +6.24%  "

$ws.Range("D9").Value = "45.48"
$ws.Range("E9").Value = "  -0.94%  "

$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").Value = "0.07471"
$ws.Range("E11").Value = "  +3.74%  "

$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").Value = "22.32"
$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("D14").Value = "6.203"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "1.793.07"
$ws.Range("E16").Value = "  +1.79%  "

$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("D18").Value = "0.06678"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "80.98"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").Value = "6.374"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").Value = "28.209.27"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").Value = "2.388"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "20.39"
$ws.Range("E26").Value = "  +3.03%  "

$ws.Range("D27").Value = "153.58"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "2.373"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").Value = "1.996.37"
$ws.Range("E29").Value = "  +1.79%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.266"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "132.24"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("D33").Value = "5.867"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").Value = "0.09389"
$ws.Range("E34").Value = "  +7.32%  "

$ws.Range("D35").Value = "0.02365"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").Value = "12.06"
$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("D37").Value = "0.6629"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").Value = "0.06225"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.168"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2154"
$ws.Range("E40").Value = "  +1.99%  "

$ws.Range("D41").Value = "1.481"
$ws.Range("E41").Value = "  +2.29%  "

$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("D43").Value = "8.052"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").Value = "13.89"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").Value = "3.861"
$ws.Range("E46").Value = "  +0.65%  "

$ws.Range("D47").Value = "0.6056"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("D48").Value = "128.36"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("D49").Value = "2.017"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").Value = "0.07080"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").Value = "1.161"
$ws.Range("E51").Value = "  -1.41%  "
